$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 gets data previously at row 3
$ws.Range("D2").Value = 44208
$ws.Range("J2").Value = 160
$ws.Range("K2").Value = 5000
$ws.Range("L2").Value = 6000
$ws.Range("M2").Value = 5500
$ws.Range("O2").Value = "Provincia de Quillota"
$ws.Range("P2").Value = 344

# Row 3 gets data previously at row 5
$ws.Range("D3").Value = 44187
$ws.Range("J3").Value = 160
$ws.Range("K3").Value = 5000
$ws.Range("L3").Value = 6000
$ws.Range("M3").Value = 5500
$ws.Range("O3").Value = "Provincia de Quillota"
$ws.Range("P3").Value = 344

# Row 4 gets data previously at row 9
$ws.Range("D4").Value = 44230
$ws.Range("J4").Value = 250
$ws.Range("K4").Value = 5000
$ws.Range("L4").Value = 6000
$ws.Range("M4").Value = 5500
$ws.Range("O4").Value = "Provincia de Quillota"
$ws.Range("P4").Value = 344

# Row 5 gets data previously at row 8
$ws.Range("D5").Value = 44210
$ws.Range("J5").Value = 340
$ws.Range("K5").Value = 5000
$ws.Range("L5").Value = 6000
$ws.Range("M5").Value = 5500
$ws.Range("O5").Value = "Provincia de Quillota"
$ws.Range("P5").Value = 344

# Row 6 gets data previously at row 7
$ws.Range("D6").Value = 44215
$ws.Range("J6").Value = 250
$ws.Range("K6").Value = 5000
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = 5500
$ws.Range("O6").Value = "Provincia de Quillota"
$ws.Range("P6").Value = 344

# Row 7 gets data previously at row 4
$ws.Range("D7").Value = 44188
$ws.Range("J7").Value = 210
$ws.Range("K7").Value = 5000
$ws.Range("L7").Value = 6000
$ws.Range("M7").Value = 5500
$ws.Range("O7").Value = "Provincia de Quillota"
$ws.Range("P7").Value = 344

# Row 8 gets data previously at row 11
$ws.Range("D8").Value = 44204
$ws.Range("J8").Value = 430
$ws.Range("K8").Value = 5000
$ws.Range("L8").Value = 6000
$ws.Range("M8").Value = 5500
$ws.Range("O8").Value = "Provincia de Quillota"
$ws.Range("P8").Value = 344

# Row 9 gets data previously at row 10
$ws.Range("D9").Value = 44292
$ws.Range("J9").Value = 90
$ws.Range("K9").Value = 6000
$ws.Range("L9").Value = 6000
$ws.Range("M9").Value = 6000
$ws.Range("O9").Value = "Región Metropolitana"
$ws.Range("P9").Value = 375

# Row 10 gets data previously at row 12
$ws.Range("D10").Value = 44231
$ws.Range("J10").Value = 250
$ws.Range("K10").Value = 5000
$ws.Range("L10").Value = 6000
$ws.Range("M10").Value = 5500
$ws.Range("O10").Value = "Provincia de Quillota"
$ws.Range("P10").Value = 344

# Row 11 gets data previously at row 2
$ws.Range("D11").Value = 44251
$ws.Range("J11").Value = 120
$ws.Range("K11").Value = 5000
$ws.Range("L11").Value = 5000
$ws.Range("M11").Value = 5000
$ws.Range("O11").Value = "Región Metropolitana"
$ws.Range("P11").Value = 312

# Row 12 gets data previously at row 6
$ws.Range("D12").Value = 44232
$ws.Range("J12").Value = 250
$ws.Range("K12").Value = 5000
$ws.Range("L12").Value = 6000
$ws.Range("M12").Value = 5500
$ws.Range("O12").Value = "Provincia de Quillota"
$ws.Range("P12").Value = 344

# Row 13 gets data previously at row 14
$ws.Range("D13").Value = 44189
$ws.Range("J13").Value = 250
$ws.Range("K13").Value = 5000
$ws.Range("L13").Value = 6000
$ws.Range("M13").Value = 5500
$ws.Range("O13").Value = "Provincia de Quillota"
$ws.Range("P13").Value = 344

# Row 14 gets data previously at row 13
$ws.Range("D14").Value = 44186
$ws.Range("J14").Value = 160
$ws.Range("K14").Value = 5000
$ws.Range("L14").Value = 6000
$ws.Range("M14").Value = 5500
$ws.Range("O14").Value = "Provincia de Quillota"
$ws.Range("P14").Value = 344

